# ModelComponentClassDiagram.pptx refactor:
#   - bump the cached "last saved" date field (datetimeFigureOut) forward a
#     week, from 10/16/18 to 10/23/18, everywhere it is rendered (slide
#     master, every slide layout, and the notes master)
#   - rename the "EnrolledClass" association-class label to "EnrolledModule"
#     on slide 1, shrinking its font so the longer word still fits the box

$p = $ppt.ActivePresentation

$newDate = "10/23/18"
$ppPlaceholderDate = 16   # PpPlaceholderType.ppPlaceholderDate
$msoPlaceholder = 14      # MsoShapeType.msoPlaceholder

function Set-DatePlaceholderText {
    param($container, [string]$text)

    for ($k = 1; $k -le $container.Shapes.Count; $k++) {
        $shp = $container.Shapes.Item($k)
        if ($shp.Type -eq $msoPlaceholder -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = $text
            break
        }
    }
}

# Slide master's date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master $newDate

# Every custom (slide) layout has its own date placeholder instance.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DatePlaceholderText $layout $newDate
}

# Notes master's date placeholder.
$notesMaster = $p.NotesMaster
Set-DatePlaceholderText $notesMaster $newDate

# Slide 1: "EnrolledClass" -> "EnrolledModule" (and shrink the font so the
# longer word keeps fitting inside the small association-class rectangle).
$slide = $p.Slides.Item(1)
$shape = $slide.Shapes.Item(46)
$tr = $shape.TextFrame.TextRange
$tr.Text = "EnrolledModule"
$tr.Font.Size = 6
